# Updated cryptos list values per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.288.44"
$ws.Range("D2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'1.831.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -0.06%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'243.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.60%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'0.6189"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.07369"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -1.12%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'0.2931"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").Value = "'23.29"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +1.04%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.07648"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'1.843.09"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.61%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'4.983"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.34%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'0.6756"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +0.28%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'82.69"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -0.26%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'0.000008976"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -2.07%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("E17").Value = "'  -0.52%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'29.280.63"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.81%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'2.091.06"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +1.00%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("D20").Value = "'237.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.60%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("E21").Value = "'  -1.26%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'0.9997"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.00%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'7.395"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +2.65%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'0.9997"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +0.13%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'158.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.33%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("D26").Value = "'0.1399"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -0.72%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").Value = "'8.561"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +0.67%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = "'17.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -1.39%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("D30").Value = "'0.05799"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +3.29%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").Value = "'1.230"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +1.99%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").Value = "'4.097"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -0.45%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("D33").Value = "'4.100"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -0.80%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").Value = "'1.866"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.28%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").Value = "'  -0.34%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("D36").Value = "'0.7228"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -2.41%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").Value = "'2.614"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.57%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").Value = "'2.859"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +3.24%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'1.225.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +1.02%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.01764"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -1.22%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("D41").Value = "'6.226"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -2.71%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'0.9080"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +1.32%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("D43").Value = "'0.9998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'2.014.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +1.97%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'101.90"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +0.43%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").Value = "'65.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.52%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").Value = "'0.5059"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.38%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("B48").Value = "'TheSandbox"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'0.4051"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.30%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Value = "'9.163"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.07%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("B50").Value = "'Algorand"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").Value = "'0.1171"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +5.58%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("B51").Value = "'BabyDogeCoin"
$ws.Range("B51").ClearFormats()
$ws.Range("C51").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").Value = "'0.00000000117"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.13%  "
$ws.Range("E51").ClearFormats()

